# Update ObjTables header attribute names from UpperCamelCase to lowerCamelCase
# and bump the embedded date stamps, per commit: "changing document, table
# attributes to lowerCamelCase".
#
# Each worksheet in this workbook corresponds to one ObjTables "table" and
# carries a single header string in cell A1 (the first sheet additionally
# carries the workbook-level "!!!ObjTables ..." header in A2).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> (cell address, new header text)
$updates = @{
    "!!Compartment" = @{
        "A1" = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 13:00:57'"
        "A2" = "!!ObjTables type='Data' id='Compartment' name='Compartment' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Compound" = @{
        "A1" = "!!ObjTables type='Data' id='Compound' name='Compound' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Definition" = @{
        "A1" = "!!ObjTables type='Data' id='Definition' name='Definition' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Enzyme" = @{
        "A1" = "!!ObjTables type='Data' id='Enzyme' name='Enzyme' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!FbcObjective" = @{
        "A1" = "!!ObjTables type='Data' id='FbcObjective' name='FbcObjective' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Gene" = @{
        "A1" = "!!ObjTables type='Data' id='Gene' name='Gene' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Layout" = @{
        "A1" = "!!ObjTables type='Data' id='Layout' name='Layout' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Measurement" = @{
        "A1" = "!!ObjTables type='Data' id='Measurement' name='Measurement' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!PbConfig" = @{
        "A1" = "!!ObjTables type='Data' id='PbConfig' name='PbConfig' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Position" = @{
        "A1" = "!!ObjTables type='Data' id='Position' name='Position' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Protein" = @{
        "A1" = "!!ObjTables type='Data' id='Protein' name='Protein' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Quantity" = @{
        "A1" = "!!ObjTables type='Data' id='Quantity' name='Quantity' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!QuantityInfo" = @{
        "A1" = "!!ObjTables type='Data' id='QuantityInfo' name='QuantityInfo' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!QuantityMatrix" = @{
        "A1" = "!!ObjTables type='Data' id='QuantityMatrix' name='QuantityMatrix' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Reaction" = @{
        "A1" = "!!ObjTables type='Data' id='Reaction' name='Reaction' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!ReactionStoichiometry" = @{
        "A1" = "!!ObjTables type='Data' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Regulator" = @{
        "A1" = "!!ObjTables type='Data' id='Regulator' name='Regulator' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Relation" = @{
        "A1" = "!!ObjTables type='Data' id='Relation' name='Relation' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!Relationship" = @{
        "A1" = "!!ObjTables type='Data' id='Relationship' name='Relationship' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrix" = @{
        "A1" = "!!ObjTables type='Data' id='SparseMatrix' name='SparseMatrix' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrixColumn" = @{
        "A1" = "!!ObjTables type='Data' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrixOrdered" = @{
        "A1" = "!!ObjTables type='Data' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrixRow" = @{
        "A1" = "!!ObjTables type='Data' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-03-09 13:00:57' objTablesVersion='0.0.8'"
    }
    "!!StoichiometricMatrix" = @{
        "A1" = "!!ObjTables type='Data' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-03-09 13:00:58' objTablesVersion='0.0.8'"
    }
    "!!rxnconContingencyList" = @{
        "A1" = "!!ObjTables type='Data' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-03-09 13:00:58' objTablesVersion='0.0.8'"
    }
    "!!rxnconReactionList" = @{
        "A1" = "!!ObjTables type='Data' id='rxnconReactionList' name='rxnconReactionList' date='2020-03-09 13:00:58' objTablesVersion='0.0.8'"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]

    # Sheets are protected (no password) - unprotect, edit, then restore
    # protection so the resulting workbook keeps the same protection state.
    $ws.Unprotect()
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
    $ws.Protect()
}
